$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.6
$summary.Range("B4").Value = -0.4
$summary.Range("B5").Value = -0.57
$summary.Range("B6").Value = 14
$summary.Range("B8").Value = 7
$summary.Range("B9").Value = 28.57

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.59999999999999
$status.Range("D4").Value = 14
$status.Range("E4").Value = -0.4
$status.Range("F4").Value = -0.4
$status.Range("G4").Value = 28.57

# --- New trade row (#14) appended to "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(15, 1).Value = 14

    # Dates get auto-detected by Excel and converted to serial numbers;
    # force the cell to Text format first so the literal string is kept.
    $ws.Cells.Item(15, 2).NumberFormat = "@"
    $ws.Cells.Item(15, 2).Value = "2026-02-17"
    $ws.Cells.Item(15, 2).ClearFormats()

    $ws.Cells.Item(15, 3).Value = "08:14:05"
    $ws.Cells.Item(15, 4).Value = "MarketMaking"
    $ws.Cells.Item(15, 5).Value = "UP"
    $ws.Cells.Item(15, 6).Value = 0.65
    $ws.Cells.Item(15, 7).Value = 0.61
    $ws.Cells.Item(15, 8).Value = "CLOSED"
    $ws.Cells.Item(15, 9).Value = -6.1538
    $ws.Cells.Item(15, 10).Value = -0.04
    $ws.Cells.Item(15, 11).Value = 99.59999999999999
    $ws.Cells.Item(15, 12).Value = 0
    $ws.Cells.Item(15, 13).Value = 0
    $ws.Cells.Item(15, 14).Value = 0.6
    $ws.Cells.Item(15, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(15, 16).Value = "early_exit"
    $ws.Cells.Item(15, 17).Value = 0.13
}
